# Generate Report for Handoff
# ----------------------------------------------------------------------------
# This script updates the localization-status workbook so that the
# "b63423d5-0118-46e2-b249-8d8ab9845102.md" file record now appears first
# (row 2) on every sheet, while the "6f3c7c99-3344-4975-854e-52c6f46f233c.md"
# record (now row 3 on every sheet) reflects a fresh handoff: status is
# "Ready for handoff", new handoff timestamps, and a populated error detail
# noting the handback file version is stale.
# ----------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$urlBase6f3c = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7b297801af86cdff3cd12880e052e7321b857f95/e2e/6f3c7c99-3344-4975-854e-52c6f46f233c.md"
$urlBaseB634 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7b297801af86cdff3cd12880e052e7321b857f95/e2e/b63423d5-0118-46e2-b249-8d8ab9845102.md"

$urlZhCn6f3c = "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/635aad4de83ebf38ae58d56277690be1fec7525f/e2e/6f3c7c99-3344-4975-854e-52c6f46f233c.md"
$urlZhCnB634 = "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/635aad4de83ebf38ae58d56277690be1fec7525f/e2e/b63423d5-0118-46e2-b249-8d8ab9845102.md"

$urlDeDe6f3c = "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/552debe71253bbaab0993597aa8079d64f0e8bb7/e2e/6f3c7c99-3344-4975-854e-52c6f46f233c.md"
$urlDeDeB634 = "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/552debe71253bbaab0993597aa8079d64f0e8bb7/e2e/b63423d5-0118-46e2-b249-8d8ab9845102.md"

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7b297801af86cdff3cd12880e052e7321b857f95/e2e/6f3c7c99-3344-4975-854e-52c6f46f233c.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d5e4e81a8d22e6b0a2aa949f0efb3c0bde2d224b/e2e/6f3c7c99-3344-4975-854e-52c6f46f233c.md."

# ============================================================================
# Sheet "Overview"
# ============================================================================
$ws1 = $wb.Worksheets.Item("Overview")

# Row 2 now carries the b63423d5 record, row 3 now carries the 6f3c7c99 record.
$ws1.Range("A2").Value = "b63423d5-0118-46e2-b249-8d8ab9845102.md"
$ws1.Range("B2").Value = "e2e\b63423d5-0118-46e2-b249-8d8ab9845102.md"

$ws1.Range("A3").Value = "6f3c7c99-3344-4975-854e-52c6f46f233c.md"
$ws1.Range("B3").Value = "e2e\6f3c7c99-3344-4975-854e-52c6f46f233c.md"

# The 6f3c7c99 record (row 3) has a new status/date.
$ws1.Range("E3").Value = "Ready for handoff"
$ws1.Range("F3").Value = "Ready for handoff"
$ws1.Range("G3").Value = "2016-08-27 12:45:40"

# Rebuild the hyperlinks for column B (target URLs are unchanged, only the
# display text moves along with the swapped rows).
$ws1.Range("A1").Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("B2"), $urlBase6f3c, "", "", "e2e\b63423d5-0118-46e2-b249-8d8ab9845102.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("B3"), $urlBaseB634, "", "", "e2e\6f3c7c99-3344-4975-854e-52c6f46f233c.md") | Out-Null

# ============================================================================
# Sheet "zh-cn"
# ============================================================================
$ws2 = $wb.Worksheets.Item("zh-cn")

# Row 2 now carries the b63423d5 record.
$ws2.Range("A2").Value = "b63423d5-0118-46e2-b249-8d8ab9845102.md"
$ws2.Range("G2").Value = "b63423d5-0118-46e2-b249-8d8ab9845102.d64bfdbfd32e5fb42d8707ac34b94cf37896e8e5.zh-cn.xlf"
$ws2.Range("I2").Value = "b63423d5-0118-46e2-b249-8d8ab9845102.md"
$ws2.Range("J2").Value = "b63423d5-0118-46e2-b249-8d8ab9845102.d64bfdbfd32e5fb42d8707ac34b94cf37896e8e5.zh-cn.xlf"

# Row 3 now carries the 6f3c7c99 record, refreshed for the new handoff.
$ws2.Range("A3").Value = "6f3c7c99-3344-4975-854e-52c6f46f233c.md"
$ws2.Range("C3").Value = "Ready for handoff"
$ws2.Range("G3").Value = "6f3c7c99-3344-4975-854e-52c6f46f233c.4a02b5ae97bded18aad9306d0ee9442eac8b5718.zh-cn.xlf"
$ws2.Range("H3").Value = "2016-08-27 12:45:36"
$ws2.Range("I3").Value = "6f3c7c99-3344-4975-854e-52c6f46f233c.md"
$ws2.Range("J3").Value = "6f3c7c99-3344-4975-854e-52c6f46f233c.4a02b5ae97bded18aad9306d0ee9442eac8b5718.zh-cn.xlf"
$ws2.Range("P3").Value = $errorDetail

# Rebuild hyperlinks for columns A and I.
$ws2.Range("A1").Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("A2"), $urlBase6f3c, "", "", "b63423d5-0118-46e2-b249-8d8ab9845102.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("I2"), $urlZhCn6f3c, "", "", "b63423d5-0118-46e2-b249-8d8ab9845102.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A3"), $urlBaseB634, "", "", "6f3c7c99-3344-4975-854e-52c6f46f233c.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("I3"), $urlZhCnB634, "", "", "6f3c7c99-3344-4975-854e-52c6f46f233c.md") | Out-Null

# Column P (Error Detail) widens to fit the new message.
$ws2.Columns.Item(16).ColumnWidth = 39.166666666666664

# ============================================================================
# Sheet "de-de"
# ============================================================================
$ws3 = $wb.Worksheets.Item("de-de")

# Row 2 now carries the b63423d5 record.
$ws3.Range("A2").Value = "b63423d5-0118-46e2-b249-8d8ab9845102.md"
$ws3.Range("G2").Value = "b63423d5-0118-46e2-b249-8d8ab9845102.d64bfdbfd32e5fb42d8707ac34b94cf37896e8e5.de-de.xlf"
$ws3.Range("I2").Value = "b63423d5-0118-46e2-b249-8d8ab9845102.md"
$ws3.Range("J2").Value = "b63423d5-0118-46e2-b249-8d8ab9845102.d64bfdbfd32e5fb42d8707ac34b94cf37896e8e5.de-de.xlf"

# Row 3 now carries the 6f3c7c99 record, refreshed for the new handoff.
$ws3.Range("A3").Value = "6f3c7c99-3344-4975-854e-52c6f46f233c.md"
$ws3.Range("C3").Value = "Ready for handoff"
$ws3.Range("G3").Value = "6f3c7c99-3344-4975-854e-52c6f46f233c.4a02b5ae97bded18aad9306d0ee9442eac8b5718.de-de.xlf"
$ws3.Range("H3").Value = "2016-08-27 12:45:40"
$ws3.Range("I3").Value = "6f3c7c99-3344-4975-854e-52c6f46f233c.md"
$ws3.Range("J3").Value = "6f3c7c99-3344-4975-854e-52c6f46f233c.4a02b5ae97bded18aad9306d0ee9442eac8b5718.de-de.xlf"
$ws3.Range("P3").Value = $errorDetail

# Rebuild hyperlinks for columns A and I.
$ws3.Range("A1").Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Range("A2"), $urlBase6f3c, "", "", "b63423d5-0118-46e2-b249-8d8ab9845102.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("I2"), $urlDeDe6f3c, "", "", "b63423d5-0118-46e2-b249-8d8ab9845102.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A3"), $urlBaseB634, "", "", "6f3c7c99-3344-4975-854e-52c6f46f233c.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("I3"), $urlDeDeB634, "", "", "6f3c7c99-3344-4975-854e-52c6f46f233c.md") | Out-Null

# Column P (Error Detail) widens to fit the new message.
$ws3.Columns.Item(16).ColumnWidth = 39.166666666666664
